$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so purely numeric-looking
# strings (e.g. "604.36", "7.00", "0.0000193") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.450.19"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "3.507.66"
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "604.36"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "143.11"
$ws.Range("E6").Value = "  -5.01%  "
$ws.Range("D7").Value = "3.505.86"
$ws.Range("E7").Value = "  -3.11%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("D10").Value = "7.69"
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("E11").Value = "  -5.81%  "
$ws.Range("D12").Value = "0.401"
$ws.Range("E12").Value = "  -3.52%  "
$ws.Range("D13").Value = "4.101.16"
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").Value = "0.0000193"
$ws.Range("E14").Value = "  -7.86%  "
$ws.Range("D15").Value = "28.49"
$ws.Range("E15").Value = "  -4.95%  "
$ws.Range("D16").Value = "3.504.67"
$ws.Range("E16").Value = "  -3.17%  "
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "66.328.42"
$ws.Range("D19").Value = "10.59"
$ws.Range("E19").Value = "  -9.22%  "
$ws.Range("D20").Value = "6.07"
$ws.Range("E20").Value = "  -4.78%  "
$ws.Range("D21").Value = "14.53"
$ws.Range("E21").Value = "  -3.77%  "
$ws.Range("D22").Value = "419.52"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "0.584"
$ws.Range("E23").Value = "  -5.70%  "
$ws.Range("D24").Value = "76.63"
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("D25").Value = "3.653.14"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "0.0000112"
$ws.Range("E27").Value = "  -9.27%  "
$ws.Range("D28").Value = "2.44"
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "8.88"
$ws.Range("E29").Value = "  -7.63%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "7.73"
$ws.Range("E30").Value = "  -7.83%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "3.515.08"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("E33").Value = "  -3.92%  "
$ws.Range("D34").Value = "24.08"
$ws.Range("E34").Value = "  -5.37%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "1.31"
$ws.Range("E36").Value = "  -10.86%  "
$ws.Range("D37").Value = "7.45"
$ws.Range("E37").Value = "  -5.60%  "
$ws.Range("D38").Value = "1.62"
$ws.Range("E38").Value = "  -5.00%  "
$ws.Range("D39").Value = "173.53"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").Value = "5.15"
$ws.Range("E40").Value = "  -8.53%  "
$ws.Range("D41").Value = "0.0801"
$ws.Range("E41").Value = "  -7.16%  "
$ws.Range("D42").Value = "4.91"
$ws.Range("E42").Value = "  -6.19%  "
$ws.Range("D43").Value = "0.848"
$ws.Range("E43").Value = "  -5.80%  "
$ws.Range("D44").Value = "45.46"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("D45").Value = "1.75"
$ws.Range("E45").Value = "  -7.70%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "2.32"
$ws.Range("E47").Value = "  -10.91%  "
$ws.Range("D48").Value = "7.00"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "22.79"
$ws.Range("E49").Value = "  -4.77%  "
$ws.Range("D50").Value = "1.10"
$ws.Range("E50").Value = "  -5.27%  "
$ws.Range("D51").Value = "0.888"
$ws.Range("E51").Value = "  -8.03%  "

# Restore the default (Normal) style so cells keep no explicit style override,
# matching the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
